# JulianSkeleWarHours.xlsx - log this week's (Week 6) hours and add its total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 was a blank pre-formatted template row; turn it into the
# "Week 6" section header, matching the style of the other week headers.
$ws.Range("A27").Value = "Week 6"
$ws.Range("A27").Font.Bold = $true

# Row 28: first logged entry for week 6.
$ws.Range("A28").Value = 42244
$ws.Range("B28").Value = 0.75
$ws.Range("C28").Value = 0.95833333333333337
$ws.Range("D28").Value = 0
$ws.Range("E28").Formula = "=MOD(C28-B28,1)*24-D28"

# Row 29: second logged entry for week 6.
$ws.Range("A29").Value = 42247
$ws.Range("B29").Value = 0.5
$ws.Range("C29").Value = 0.95833333333333337
$ws.Range("D29").Value = 2
$ws.Range("E29").Formula = "=MOD(C29-B29,1)*24-D29"

# Row 30: new "Total" row summing week 6's hours.
$ws.Range("D30").Value = "Total"
$ws.Range("D30").Font.Bold = $true
$ws.Range("E30").Formula = "=SUM(E28:E29)"
$ws.Range("E30").Font.Bold = $true
$ws.Range("E30").NumberFormat = "0.00"

# Move the selection to where the user left off after typing the new total.
$ws.Range("F31").Select() | Out-Null
